# Update countries & provincias Spain
# Applies updated COVID case counters and re-sorts a few tied rows
# (Hong Kong/Nicaragua, Lituania/Guadalupe, Santa Lucia/Timor Oriental)
# as well as the "last updated" timestamp in the title cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 11:29"

# --- Rows whose data updated, no reordering needed ---

# Row 25: Alemania
$ws.Range("B25").Value = 286420
$ws.Range("C25").Value = 82
$ws.Range("D25").Value = 250800
$ws.Range("E25").Value = 26086
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 9534

# Row 26: Indonesia
$ws.Range("B26").Value = 278722
$ws.Range("C26").Value = 3509
$ws.Range("D26").Value = 206870
$ws.Range("E26").Value = 61379
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 87
$ws.Range("H26").Value = 10473

# Row 46: Polonia
$ws.Range("B46").Value = 88636
$ws.Range("C46").Value = 1306
$ws.Range("D46").Value = 68420
$ws.Range("E46").Value = 17769
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 15
$ws.Range("H46").Value = 2447

# Row 67: Austria
$ws.Range("B67").Value = 43432
$ws.Range("C67").Value = 556
$ws.Range("D67").Value = 34052
$ws.Range("E67").Value = 8590
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 790

# Row 106: Eslovaquia
$ws.Range("B106").Value = 9343
$ws.Range("C106").Value = 265
$ws.Range("D106").Value = 4213
$ws.Range("E106").Value = 5086
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 44

# Row 142: Sri Lanka
$ws.Range("B142").Value = 3360
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 3210
$ws.Range("E142").Value = 137
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 13

# --- Rows 122-124: Eslovenia stays; Hong Kong and Nicaragua swap order
#     (Hong Kong's totals overtook Nicaragua's) ---

# Row 122: Eslovenia (data updated)
$ws.Range("A122").Value = "Eslovenia"
$ws.Range("B122").Value = 5388
$ws.Range("C122").Value = 39
$ws.Range("D122").Value = 3600
$ws.Range("E122").Value = 1639
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 149

# Row 123: now Hong Kong (was Nicaragua), with Hong Kong's updated data
$ws.Range("A123").Value = "Hong Kong"
$ws.Range("B123").Value = 5076
$ws.Range("C123").Value = 10
$ws.Range("D123").Value = 4790
$ws.Range("E123").Value = 181
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 105

# Row 124: now Nicaragua (was Hong Kong), keeping Nicaragua's original data
$ws.Range("A124").Value = "Nicaragua"
$ws.Range("B124").Value = 5073
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 2913
$ws.Range("E124").Value = 2011
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 149

# --- Rows 131-132: Lituania and Guadalupe swap order
#     (Lituania's totals overtook Guadalupe's) ---

# Row 131: now Lituania (was Guadalupe), with Lituania's updated data
$ws.Range("A131").Value = "Lituania"
$ws.Range("B131").Value = 4490
$ws.Range("C131").Value = 105
$ws.Range("D131").Value = 2327
$ws.Range("E131").Value = 2071
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 92

# Row 132: now Guadalupe (was Lituania), keeping Guadalupe's original data
$ws.Range("A132").Value = "Guadalupe"
$ws.Range("B132").Value = 4487
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 2199
$ws.Range("E132").Value = 2246
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 42

# --- Rows 207-208: Santa Lucia and Timor Oriental swap order
#     (tied data, alphabetical tie-break) ---

# Row 207: now Santa Lucia (was Timor Oriental)
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("B207").Value = 27
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 27
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

# Row 208: now Timor Oriental (was Santa Lucia)
$ws.Range("A208").Value = "Timor Oriental"
$ws.Range("B208").Value = 27
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 27
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
